$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing E9 value (seasonal stat correction)
$ws.Range("E9").Value = 838212

# Copy the formatting of the A9 cell (bold, bordered, centered style) onto the
# new A10 cell before setting its value, so the new row's season-index column
# matches the existing style.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new row 10 data
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "M2_09 Dryad 2020"
$ws.Range("C10").Value = 9678
$ws.Range("D10").Value = 10725
$ws.Range("E10").Value = 855528
$ws.Range("F10").Value = 9946
$ws.Range("G10").Value = 10046
$ws.Range("H10").Value = 10183
